$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "26.015.51"
$ws.Range("E2").Value = "  -0.04%  "

# Row 3
$ws.Range("D3").Value = "1.635.05"
$ws.Range("E3").Value = "  -0.52%  "

# Row 4
$ws.Range("E4").Value = "  -0.02%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "214.22"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.97%  "

# Row 6
$ws.Range("E6").Value = "  -0.66%  "

# Row 7
$ws.Range("E7").Value = "  -0.02%  "

# Row 8
$ws.Range("E8").Value = "  -2.19%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.0625"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -2.22%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "18.54"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -5.60%  "

# Row 11
$ws.Range("E11").Value = "  -0.44%  "

# Row 12
$ws.Range("D12").Value = "1.862.02"

# Row 13
$ws.Range("E13").Value = "  -1.87%  "

# Row 14
$ws.Range("D14").Value = "1.633.35"
$ws.Range("E14").Value = "  -1.58%  "

# Row 15
$ws.Range("E15").Value = "  -2.84%  "

# Row 16
$ws.Range("E16").Value = "  -2.53%  "

# Row 17
$ws.Range("D17").Value = "26.010.45"
$ws.Range("E17").Value = "  +0.30%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "61.90"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.80%  "

# Row 19
$ws.Range("E19").Value = "  -0.02%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "190.88"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.18%  "

# Row 21
$ws.Range("E21").Value = "  -2.06%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.59"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -3.58%  "

# Row 23
$ws.Range("E23").Value = "  -1.85%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.135"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.88%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "143.24"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.15%  "

# Row 26
$ws.Range("E26").Value = "  +0.03%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.77"
$ws.Range("D27").Style = "Normal"

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "6.76"
$ws.Range("D28").Style = "Normal"

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.22"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -2.30%  "

# Row 30
$ws.Range("E30").Value = "  -1.40%  "

# Row 31
$ws.Range("E31").Value = "  -3.50%  "

# Row 32
$ws.Range("E32").Value = "  -2.77%  "

# Row 33
$ws.Range("E33").Value = "  -3.99%  "

# Row 34
$ws.Range("B34").Value = "LidoDAOToken"
$ws.Range("C34").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.51"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -2.04%  "

# Row 35
$ws.Range("B35").Value = "HuobiToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.43"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.67%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.872"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -3.63%  "

# Row 37
$ws.Range("D37").Value = "1.132.83"
$ws.Range("E37").Value = "  -0.24%  "

# Row 38
$ws.Range("B38").Value = "ImmutableX"
$ws.Range("C38").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.527"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -3.12%  "

# Row 39
$ws.Range("B39").Value = "MXToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.43"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.64%  "

# Row 40
$ws.Range("E40").Value = "  -1.52%  "

# Row 41
$ws.Range("E41").Value = "  -1.14%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.786"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.53%  "

# Row 43
$ws.Range("E43").Value = "  -4.21%  "

# Row 44
$ws.Range("D44").Value = "1.772.16"

# Row 45
$ws.Range("E45").Value = "  -1.07%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "55.31"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -2.68%  "

# Row 47
$ws.Range("E47").Value = "  -0.69%  "

# Row 48
$ws.Range("E48").Value = "  +1.43%  "

# Row 49
$ws.Range("E49").Value = "  -0.55%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.52"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.71%  "

# Row 51
$ws.Range("E51").Value = "  -0.19%  "
